$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E14: plain value change 637985.22 -> 537985.22
$ws.Range("E14").Value = 537985.22

# E18: was a plain value, becomes formula SUM(E12:E17)
$ws.Range("E18").Formula = "=SUM(E12:E17)"

# E21: was a plain value, becomes formula SUM(E18:E20)
$ws.Range("E21").Formula = "=SUM(E18:E20)"

# E26: plain value change 101863404 -> 1018613404
$ws.Range("E26").Value = 1018613404

$wb.Save()
